# "Generate Report for Handoff"
#
# Localization status has moved from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated timestamps have advanced to
# reflect the new handoff generation run. Update the Overview sheet and
# each per-locale sheet (zh-cn, de-de) accordingly, and shrink the
# now-shorter Status column(s) to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# E2/F2 hold the per-locale status ("zh-cn" / "de-de" columns), G2 holds the
# latest handoff xliff-generation timestamp for the row.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 17:40:59"

# --- zh-cn sheet ---
# C2 is Status, H2 is Latest Handoff Datetime.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 17:40:54"

# --- de-de sheet ---
# C2 is Status, H2 is Latest Handback DateTime.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 17:40:59"

# The Status column is now narrower text ("Ready for handoff" vs.
# "Handed back: in sync with en-US"), so the authoring tool shrank the
# corresponding columns from ~30 characters wide to ~17.
$overview.Range("E1").ColumnWidth = 16.4
$overview.Range("F1").ColumnWidth = 16.4
$zhcn.Range("C1").ColumnWidth = 16.4
$dede.Range("C1").ColumnWidth = 16.4
